# Scheduled runner update: refresh market-board price snapshots and
# recompute dependent profit columns across the Leve sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5049.375
$ws.Range("I40").Value = 2449.5
$ws.Range("K40").Value = 2449.5
$ws.Range("M40").Value = -2274.5

$ws.Range("I61").Value = 500
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1500
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = ""
$ws.Range("N61").Value = -1328

$ws.Range("H62").Value = 2420.375
$ws.Range("I62").Value = 2194.7144
$ws.Range("K62").Value = 2194.7144
$ws.Range("M62").Value = -1570.7144

$ws.Range("H65").Value = 2420.375
$ws.Range("I65").Value = 2194.7144
$ws.Range("K65").Value = 10973.572
$ws.Range("M65").Value = -7853.572

$ws.Range("H69").Value = 1200
$ws.Range("I69").Value = 1200
$ws.Range("K69").Value = 3600
$ws.Range("M69").Value = -2726

$ws.Range("H72").Value = 1200
$ws.Range("I72").Value = 1200
$ws.Range("K72").Value = 10800
$ws.Range("M72").Value = -6432

$ws.Range("H74").Value = 6738.7144
$ws.Range("I74").Value = 5512
$ws.Range("K74").Value = 5512
$ws.Range("M74").Value = -4576

$ws.Range("H77").Value = 6738.7144
$ws.Range("I77").Value = 5512
$ws.Range("K77").Value = 27560
$ws.Range("M77").Value = -22880

$ws.Range("H135").Value = 705
$ws.Range("I135").Value = 679.8333
$ws.Range("K135").Value = 6118.4997
$ws.Range("M135").Value = -3583.4997

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1151
$ws.Range("I2").Value = 1151.8
$ws.Range("J2").Value = 1149
$ws.Range("K2").Value = 1151.8
$ws.Range("L2").Value = 1149
$ws.Range("M2").Value = -1038.8
$ws.Range("N2").Value = -1375

$ws.Range("H25").Value = 7671.6665
$ws.Range("J25").Value = 18999
$ws.Range("L25").Value = 18999
$ws.Range("N25").Value = -19803

$ws.Range("H46").Value = 12508.571
$ws.Range("I46").Value = 28068.5
$ws.Range("J46").Value = 6284.6
$ws.Range("K46").Value = 28068.5
$ws.Range("L46").Value = 6284.6
$ws.Range("M46").Value = -27749.5
$ws.Range("N46").Value = -6922.6

$ws.Range("H116").Value = 1151
$ws.Range("I116").Value = 1151.8
$ws.Range("J116").Value = 1149
$ws.Range("K116").Value = 1151.8
$ws.Range("L116").Value = 1149
$ws.Range("M116").Value = 1142.2
$ws.Range("N116").Value = -5737

$ws.Range("H124").Value = 70000
$ws.Range("J124").Value = 70000
$ws.Range("L124").Value = 70000
$ws.Range("N124").Value = -79820

$ws.Range("H125").Value = 75070.836
$ws.Range("J125").Value = 75070.836
$ws.Range("L125").Value = 75070.836
$ws.Range("N125").Value = -84910.836

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1151
$ws.Range("I3").Value = 1151.8
$ws.Range("J3").Value = 1149
$ws.Range("K3").Value = 1151.8
$ws.Range("L3").Value = 1149
$ws.Range("M3").Value = -1037.8
$ws.Range("N3").Value = -1377

$ws.Range("H107").Value = 350
$ws.Range("I107").Value = 350
$ws.Range("K107").Value = 350
$ws.Range("M107").Value = 1570

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 825.375
$ws.Range("I22").Value = 771.8570999999999
$ws.Range("K22").Value = 771.8570999999999
$ws.Range("M22").Value = -421.8570999999999

$ws.Range("H134").Value = 8930764
$ws.Range("I134").Value = 9525815
$ws.Range("J134").Value = 4999
$ws.Range("K134").Value = 28577445
$ws.Range("L134").Value = 14997
$ws.Range("M134").Value = -28574910
$ws.Range("N134").Value = -20067

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 114.29412
$ws.Range("I2").Value = 81.44444
$ws.Range("J2").Value = 151.25
$ws.Range("K2").Value = 488.66664
$ws.Range("L2").Value = 907.5
$ws.Range("M2").Value = -375.66664
$ws.Range("N2").Value = -1133.5

$ws.Range("H11").Value = 2014.3334
$ws.Range("I11").Value = 2014.3334
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 6043.0002
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = ""
$ws.Range("N11").Value = -5903.0002

$ws.Range("H17").Value = 68.2
$ws.Range("I17").Value = 68.2
$ws.Range("K17").Value = 204.6
$ws.Range("M17").Value = -35.60000000000002

$ws.Range("H68").Value = 2944.2593
$ws.Range("J68").Value = 3122.7727
$ws.Range("L68").Value = 9368.3181
$ws.Range("N68").Value = -10990.3181

$ws.Range("H71").Value = 2944.2593
$ws.Range("J71").Value = 3122.7727
$ws.Range("L71").Value = 28104.9543
$ws.Range("N71").Value = -36216.9543

$ws.Range("H113").Value = 524.3570999999999
$ws.Range("I113").Value = 244.4
$ws.Range("J113").Value = 679.8889
$ws.Range("K113").Value = 733.2
$ws.Range("L113").Value = 2039.6667
$ws.Range("M113").Value = 1436.8
$ws.Range("N113").Value = -6379.6667

$ws.Range("H122").Value = 814.0769
$ws.Range("I122").Value = 603.8570999999999
$ws.Range("K122").Value = 5434.7139
$ws.Range("M122").Value = -2984.7139

$ws.Range("H131").Value = 2103.0557
$ws.Range("J131").Value = 2170.5557
$ws.Range("L131").Value = 6511.6671
$ws.Range("N131").Value = -16591.6671

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 4666.3335
$ws.Range("I55").Value = 1000
$ws.Range("J55").Value = 6499.5
$ws.Range("K55").Value = 1000
$ws.Range("L55").Value = 6499.5
$ws.Range("M55").Value = -827
$ws.Range("N55").Value = -6845.5

$ws.Range("H61").Value = 3410.3572
$ws.Range("I61").Value = 3249.6667
$ws.Range("J61").Value = 4374.5
$ws.Range("K61").Value = 3249.6667
$ws.Range("L61").Value = 4374.5
$ws.Range("M61").Value = -3047.6667
$ws.Range("N61").Value = -4778.5

$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = ""
$ws.Range("N62").Value = 0

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = ""
$ws.Range("N65").Value = 0

$ws.Range("H100").Value = 2902.1428
$ws.Range("I100").Value = 2902.1428
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 2902.1428
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = ""
$ws.Range("N100").Value = -2361.1428

$ws.Range("H113").Value = 3410.3572
$ws.Range("I113").Value = 3249.6667
$ws.Range("J113").Value = 4374.5
$ws.Range("K113").Value = 3249.6667
$ws.Range("L113").Value = 4374.5
$ws.Range("M113").Value = -1079.6667
$ws.Range("N113").Value = -8714.5

$ws.Range("H135").Value = 70000
$ws.Range("K135").Value = 70000
$ws.Range("L135").Value = 70000
$ws.Range("N135").Value = -80140

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 23597
$ws.Range("I74").Value = 20569
$ws.Range("J74").Value = 25111
$ws.Range("K74").Value = 20569
$ws.Range("L74").Value = 25111
$ws.Range("M74").Value = -19633
$ws.Range("N74").Value = -26983

$ws.Range("H77").Value = 23597
$ws.Range("I77").Value = 20569
$ws.Range("J77").Value = 25111
$ws.Range("K77").Value = 61707
$ws.Range("L77").Value = 75333
$ws.Range("M77").Value = -57027
$ws.Range("N77").Value = -84693

$ws.Range("H81").Value = 2140
$ws.Range("I81").Value = 2188.6667
$ws.Range("J81").Value = 1848
$ws.Range("K81").Value = 4377.3334
$ws.Range("L81").Value = 3696
$ws.Range("M81").Value = -3316.3334
$ws.Range("N81").Value = -5818

$ws.Range("H84").Value = 2140
$ws.Range("I84").Value = 2188.6667
$ws.Range("J84").Value = 1848
$ws.Range("K84").Value = 21886.667
$ws.Range("L84").Value = 18480
$ws.Range("M84").Value = -16582.667
$ws.Range("N84").Value = -29088

$ws.Range("H122").Value = 2859.652
$ws.Range("I122").Value = 2872.5
$ws.Range("K122").Value = 8617.5
$ws.Range("M122").Value = -6167.5

$ws.Range("H136").Value = 4764.391
$ws.Range("I136").Value = 5357.1577
$ws.Range("J136").Value = 1948.75
$ws.Range("K136").Value = 16071.4731
$ws.Range("L136").Value = 5846.25
$ws.Range("M136").Value = -13521.4731
$ws.Range("N136").Value = -10946.25
